$wb = $excel.ActiveWorkbook

# --- YDS sheet: append Wild Card play-by-play yardage logs ---
$ws = $wb.Worksheets.Item("YDS")
$ws.Range("B2").Value = "31 -1 2 8 0 13 0 0 0 2 2 2 4 1 15 3 -2 6 0 1 2 4 1 0 4 5 -2 2 1 -2 4 -2 1 1 -4 11 5 2 -4 3 13 1 3 1 3 2 3 5 3 0 4 2 6 4 2 5 2 14 7 4 1 19 3 2 0 1 4 -2 1 6 0 1 27 8 1 4 2 1 5 3 -1 5 4 0 2 2 10 1 2 0 6 2 1 3 0 3 11 3 8 6 5 0 5 -1 5 5 4 2 7 1 1 0 1 1 8 1 4 6 2 6 4 3 9 3 6 18 6 2 0 1 3 1 5 3 1 3 0 1 5 1 2 3 12 9 4 3 9 4 8 -3 6 11 11 4 1 5 1 2 7 7 5 1 3 1 1 11 6 13 3 0 13 3 7 2 1 14 4 1 0 3 4 21 2 3 3 5 3 3 1 -1 12 2 -1 9 2 1 3 8 3 6 4 5 1 7 4 1 4 3 2 4 5 0 18 -4 7 2 3 2 2 4 1 3 0 9 4 1 2 3 0 10 5 5 4 -4 6 6 3 22 2 4 4 4 2 4 5 21 -1 0 4 4 3 7 1 3 0 17 7 3 15 -1 6 2 1 1 10 1 -3 5 1 4 3 4 1 5 3 2 -2 1 3 19 7 6 3 3 6 3 5 8 6 3 3 0 1 6 2 9 4 1 -3 3 3 7 4 6 5 3 3 -2 0 -4 8 7 3 8 19 4 8 5 5 4 8 -4 3 0 4 4 4 10 6 0 4 13 3 5 6 2 6 4 6 9 2 9 0 0 2 -4 1 -1 3 5 4 4 8 3 2 1 1 2 5 2 8 0 0 -3 14 -2 3 1 1 -2 0 23 1 2 7 5 3 15 2 1 0 11 7 5 8 -4 28 18 0 7 -1 7 10 5 3 1 2 2 14 35 2 20 -1 0 8 7 5"
$ws.Range("C2").Value = "3 18 2 2 3 -3 0 4 1 16 3 2 4 11 6 -1 15 11 4 2 2 5 5 8 1 0 28 2 4 2 3 2 -1 14 -2 1 1 14 -3 2 3 1 4 0 3 6 5 2 2 3 24 4 1 6 0 1 14 10 1 3 5 3 19 4 4 0 -1 0 1 15 7 -1 18 -1 0 3 5 21 5 3 -3 8 2 11 6 4 8 0 2 2 4 1 7 6 9 1 6 11 20 4 3 -2 6 1 1 3 0 5 11 14 0 0 6 8 9 3 3 9 7 0 4 8 8 1 0 0 6 2 0 7 1 6 6 5 0 1 0 1 1 10 1 9 3 11 3 4 0 3 0 5 8 4 0 30 4 2 4 4 1 1 0 6 3 5 9 7 1 5 7 5 2 7 1 4 3 4 22 5 8 5 0 6 0 0 6 4 1 16 4 2 1 5 6 4 16 5 5 6 2 -2 2 6 1 6 20 5 2 3 5 2 0 5 7 5 2 3 8 -3 4 11 12 2 5 0 5 1 3 21 1 3 2 0 6 -1 5 3 2 6 -1 4 4 4 3 2 6 4 14 4 5 1 3 -2 -3 1 -2 10 -1 11 2 5 4 8 -2 18 5 2 6 4 -5 2 1 19 4 0 1 3 7 20 1 0 1 4 2 5 8 1 3 5 2 1 4 3 2 3 8 3 1 12 -5 1 1 22 1 3 1 16 2 6 2 2 4 1 9 5 -4 -3 13 4 1 0 -6 7 2 3 8 1 6 0 4 3 2 1 2 7 5 7 0 3 1 0 0 1 4 9 6 4 2 10 4 2 -2 2 4 2 5 2 4 -7 -2 8 1 4 4 -1 2 4 24 2 4 1 1 6 4 3 12 2 6 2 5 0 4 5 1 1 1 3 1 3 -4 1 -3 1 3 0 -4 5 3 -3 4 6 4 18 24 3 0 1 1 14 2 0 4 9 -1 2 8 5 7 1 4 3 1 3 1 9 2 2 14 1 2 14 11 8 6 -1 2 -3 6 0 7 -1 4 3 4 3 5 2 5 7 2 4 -1 9 1 10 5 4 2 2 2 3 0 -2 1"
$ws.Range("B3").Value = "9 13 6 5 6 12 9 24 8 15 13 7 4 12 16 2 3 4 21 4 5 16 9 7 9 37 10 10 20 18 11 27 32 31 9 3 11 10 17 13 10 25 6 4 3 8 9 14 18 9 23 30 9 13 2 13 16 7 4 10 61 25 24 7 23 6 10 25 8 18 8 5 23 12 21 1 4 15 10 12 23 7 23 16 3 16 32 34 4 3 12 5 21 7 0 15 8 8 5 10 9 18 3 51 5 5 4 3 12 8 29 9 29 4 4 9 6 17 5 6 18 3 2 4 6 3 12 7 10 9 48 2 3 25 11 16 6 8 8 10 33 31 29 3 40 8 51 12 4 43 -1 16 2 8 13 0 6 18 18 8 10 10 6 20 7 6 5 4 21 7 9 10 24 5 -1 29 3 6 6 20 6 2 8 3 7 8 10 20 7 2 15 9 5 9 20 24 22 11 2 3 12 3 19 9 17 8 3 6 5 1 4 3 27 3 6 7 5 22 37 24 40 11 5 5 7 9 6 12 3 15 -2 -2 25 6 8 6 21 4 1 9 2 25 31 19 7 8 20 9 8 8 12 56 5 12 12 21 17 13 31 3 11 -1 4 35 10 -5 16 54 9 9 4 30 7 8 1 8 9 3 10 5 3 9 7 5 0 18 10 -5 5 2 15 9 34 7 3 28 14 13 16 2 11 7 7 8 9 16 14 2 7 19 7 2 5 7 13 10 6 4 6 4 19 19 6 4 -1 12 14 11 17 12 7 5 5 10 5 13 8 6 16 14 16 12 6 7 3 12 9 12 15 8 5 16 18 6 4 0 15 10 20 -5 5 7 13 12 12 18 6 3 28 14 5 7 3 11 -2 11 4 26 12 5 4 5 19 13 7 42 2 5 9 4 14 24 8 44 0 4 12 12 18 4 4 30 6 2 -1 8 0 7 17 -5 11 17 6 8 9 13 7 6 7 19 4 11 14 12 6 9 16 7 5 16 3 3 13 5 13 7 15 23 10"
$ws.Range("C3").Value = "2 8 0 21 6 6 29 29 10 13 4 10 6 4 13 6 49 9 10 6 5 2 2 9 7 9 17 6 13 5 41 5 14 3 8 16 52 25 4 5 9 5 4 12 0 11 8 5 -1 12 8 4 10 9 7 4 13 4 4 2 7 8 9 2 2 3 5 9 6 7 4 10 3 7 1 5 4 22 7 7 4 3 5 5 11 34 0 10 2 9 4 7 13 9 14 9 15 4 4 18 7 2 2 15 8 12 8 9 7 13 14 7 14 23 5 1 3 6 11 7 5 7 10 14 3 14 11 6 10 5 7 26 13 12 12 12 26 6 18 16 6 4 4 1 6 24 13 8 -3 20 27 25 15 20 7 5 6 17 10 8 7 6 9 5 30 6 12 7 -7 14 3 14 2 9 1 -1 6 6 7 5 11 22 19 8 8 7 13 12 -1 16 27 3 1 0 9 7 11 7 13 17 2 10 16 32 3 7 3 38 8 30 22 3 6 3 2 7 14 9 14 3 7 6 4 17 7 3 7 9 9 12 6 1 6 15 17 10 10 9 4 51 5 13 5 7 7 15 7 41 12 7 6 5 4 32 32 11 5 3 17 10 2 2 4 23 -1 13 6 7 4 7 6 5 6 5 6 0 10 35 10 3 4 10 9 10 12 6 6 44 23 10 19 19 28 7 9 3 1 7 38 15 5 -3 2 8 16 1 5 10 6 6 2 4 3 4 4 9 13 6 6 13 5 25 11 3 6 5 2 7 1 40 4 14 6 8 20 5 13 4 8 8 5 17 7 6 8 4 12 7 45 5 3 12 6 0 3 3 5 5 14 19 6 10 6 6 13 14 7 42 2 3 6 5 9 18 1 6 14 4 14 10 17 12 6 11 17 9 21 7 -1 28 5 5 29 5 10 5 9 7 9 5 7 3 19 11 10 6 7"

# --- OFF sheet: season totals updated with Wild Card game ---
$ws = $wb.Worksheets.Item("OFF")
$ws.Range("C2").Value = 194
$ws.Range("D2").Value = 20
$ws.Range("E2").Value = 9
$ws.Range("F2").Value = 55
$ws.Range("G2").Value = 68
$ws.Range("H2").Value = 7
$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 31
$ws.Range("L2").Value = 327
$ws.Range("M2").Value = 225
$ws.Range("N2").Value = 21
$ws.Range("O2").Value = 21
$ws.Range("P2").Value = 13
$ws.Range("Q2").Value = 580
$ws.Range("B3").Value = 11
$ws.Range("C3").Value = 238
$ws.Range("E3").Value = 40
$ws.Range("F3").Value = 133
$ws.Range("G3").Value = 52
$ws.Range("H3").Value = 32
$ws.Range("I3").Value = 62
$ws.Range("J3").Value = 70
$ws.Range("L3").Value = 338
$ws.Range("M3").Value = 221
$ws.Range("N3").Value = 20
$ws.Range("Q3").Value = 573

# --- DEF sheet: season totals updated with Wild Card game ---
$ws = $wb.Worksheets.Item("DEF")
$ws.Range("C2").Value = 223
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = 70
$ws.Range("G2").Value = 73
$ws.Range("J2").Value = 39
$ws.Range("L2").Value = 327
$ws.Range("M2").Value = 209
$ws.Range("N2").Value = 19
$ws.Range("O2").Value = 30
$ws.Range("P2").Value = 22
$ws.Range("Q2").Value = 611
$ws.Range("B3").Value = 12
$ws.Range("C3").Value = 201
$ws.Range("E3").Value = 35
$ws.Range("F3").Value = 128
$ws.Range("G3").Value = 43
$ws.Range("H3").Value = 29
$ws.Range("I3").Value = 85
$ws.Range("J3").Value = 58
$ws.Range("L3").Value = 309
$ws.Range("M3").Value = 214
$ws.Range("N3").Value = 18
$ws.Range("Q3").Value = 562

# --- ST sheet: season totals + punt logs updated with Wild Card game ---
$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value = 90
$ws.Range("D2").Value = 65
$ws.Range("F2").Value = 141
$ws.Range("G2").Value = 134
$ws.Range("H2").Value = 3
$ws.Range("J2").Value = 72
$ws.Range("K2").Value = 69
$ws.Range("L2").Value = 41
$ws.Range("M2").Value = 31
$ws.Range("N2").Value = 16
$ws.Range("O2").Value = 13
$ws.Range("B3").Value = 54
$ws.Range("B4").Value = "61 61 69 66 56 62 61 61 69 63 62 65 60 66 58 65 58 65 63 61 62 59 67 65 63 66 62 58 58 57 66 59 54 57 63 64"
$ws.Range("D3").Value = "35 63 67 64 37 55 40 50 45 55 49 45 67 54 45 49 56 52 68 39 70 47 66 71 54 37 59 35 42 58 40 61 42 44 63 37 47 54 48 68 32 53 36 54 42 67 38 55 44 60 46 42 49 41 45 38 47 44 60 46 33 48 54 58 41"
$ws.Range("B5").Value = "28 18 28 18 15 24 11 24 32 13 13 18 12 23 26 23 18 18 20 14 25 32 18 0 11 24 24 18 20 0 24 23 14 27 23 25"
$ws.Range("D4").Value = "0 22 15 30 0 15 0 0 15 8 0 7 12 11 0 0 9 0 17 0 21 0 0 0 0 0 6 0 12 14 0 37 0 9 0 0 0 0 15 2 0 0 0 14 0 16 0 13 0 0 11 20 0 0 0 0 0 0 0 4 0 0 11 14 12"
$ws.Range("B6").Value = "13 20 22 15 25 15 31 22 15 15 20 14 26 32 22 20 23 19 18 19 12 25 15 19 24 28 16 29 14 23 23 34 9 18 22 0 35 16 22"
$ws.Range("D5").Value = "3 0 0 7 0 10 0 16 7 0 2 0 0 17 0 5 6 0 0 0 0 9 0 5 0 0 0 0 15 0 0 5 0 12 11 21 0 0 0 0 9 7 6 0 0 0 6 29 0 0 0 2 0 0 6 0 0 7 41 0 0 16 10 9 13 0"

# --- TURNS sheet: turnover counts updated with Wild Card game ---
$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 8
$ws.Range("B3").Value = 8
$ws.Range("D3").Value = 9
$ws.Range("E3").Value = 9

# --- PEN sheet: penalty counts updated with Wild Card game ---
$ws = $wb.Worksheets.Item("PEN")
$ws.Range("B2").Value = 21
$ws.Range("B3").Value = 28
$ws.Range("D3").Value = 9
$ws.Range("B5").Value = 2
